# Update "想去人数" (want-to-go count) figures in column F across the
# relevant worksheets, per the refreshed crawl snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibitions)
$wsExhibit.Range("F2").Value  = 587
$wsExhibit.Range("F3").Value  = 265
$wsExhibit.Range("F6").Value  = 387
$wsExhibit.Range("F9").Value  = 244
$wsExhibit.Range("F10").Value = 231
$wsExhibit.Range("F11").Value = 6100
$wsExhibit.Range("F13").Value = 53
$wsExhibit.Range("F14").Value = 508
$wsExhibit.Range("F22").Value = 165
$wsExhibit.Range("F23").Value = 98
$wsExhibit.Range("F24").Value = 320
$wsExhibit.Range("F25").Value = 1027
$wsExhibit.Range("F27").Value = 1850
$wsExhibit.Range("F28").Value = 504

# 本地生活 (Local Life)
$wsLocal.Range("F2").Value = 265

# 全部类型 (All Types)
$wsAll.Range("F2").Value  = 265
$wsAll.Range("F3").Value  = 587
$wsAll.Range("F4").Value  = 265
$wsAll.Range("F8").Value  = 387
$wsAll.Range("F11").Value = 244
$wsAll.Range("F12").Value = 231
$wsAll.Range("F13").Value = 6100
$wsAll.Range("F15").Value = 53
$wsAll.Range("F17").Value = 508
$wsAll.Range("F32").Value = 165
$wsAll.Range("F33").Value = 98
$wsAll.Range("F34").Value = 320
$wsAll.Range("F35").Value = 1027
$wsAll.Range("F37").Value = 1850
$wsAll.Range("F38").Value = 504
